$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1139.1111
$ws.Range("I137").Value = 1300.4445
$ws.Range("J137").Value = 977.7778
$ws.Range("K137").Value = 3901.3335
$ws.Range("L137").Value = 2933.3334
$ws.Range("M137").Value = -1351.3335
$ws.Range("N137").Value = -8033.3334
$ws.Range("H138").Value = 2396.45
$ws.Range("I138").Value = 1405.3829
$ws.Range("J138").Value = 3275.3208
$ws.Range("K138").Value = 4216.1487
$ws.Range("L138").Value = 9825.9624
$ws.Range("M138").Value = 923.8513000000003
$ws.Range("N138").Value = -20105.9624
$ws.Range("H141").Value = 1454.0555
$ws.Range("I141").Value = 1454.0555
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4362.166499999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 817.8335000000006
$ws.Range("N141").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4809.42
$ws.Range("I32").Value = 4707.495
$ws.Range("J32").Value = 14900
$ws.Range("K32").Value = 4707.495
$ws.Range("L32").Value = 14900
$ws.Range("M32").Value = -4420.495
$ws.Range("N32").Value = -15474
$ws.Range("H132").Value = 1686.8909
$ws.Range("I132").Value = 1275.34
$ws.Range("J132").Value = 5802.4
$ws.Range("K132").Value = 3826.02
$ws.Range("L132").Value = 17407.2
$ws.Range("M132").Value = -1296.02
$ws.Range("N132").Value = -22467.2
$ws.Range("H137").Value = 77321.07000000001
$ws.Range("I137").Value = 40000
$ws.Range("J137").Value = 80191.92
$ws.Range("K137").Value = 40000
$ws.Range("L137").Value = 80191.92
$ws.Range("M137").Value = -34900
$ws.Range("N137").Value = -90391.92

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 880.5294
$ws.Range("I99").Value = 804.46155
$ws.Range("J99").Value = 1127.75
$ws.Range("K99").Value = 804.46155
$ws.Range("L99").Value = 1127.75
$ws.Range("M99").Value = 693.53845
$ws.Range("N99").Value = -4123.75
$ws.Range("H105").Value = 3333.383
$ws.Range("I105").Value = 1894.1666
$ws.Range("J105").Value = 3826.8286
$ws.Range("K105").Value = 1894.1666
$ws.Range("L105").Value = 3826.8286
$ws.Range("M105").Value = -147.1666
$ws.Range("N105").Value = -7320.8286

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 1050
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 1050
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -1750
$ws.Range("H31").Value = 2905.8032
$ws.Range("I31").Value = 2235.4119
$ws.Range("J31").Value = 3750
$ws.Range("K31").Value = 2235.4119
$ws.Range("L31").Value = 3750
$ws.Range("M31").Value = -1940.4119
$ws.Range("N31").Value = -4340
$ws.Range("H34").Value = 2905.8032
$ws.Range("I34").Value = 2235.4119
$ws.Range("J34").Value = 3750
$ws.Range("K34").Value = 2235.4119
$ws.Range("L34").Value = 3750
$ws.Range("M34").Value = -2033.4119
$ws.Range("N34").Value = -4154

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = -431
$ws.Range("N17").Value = -938
$ws.Range("H34").Value = 779.53845
$ws.Range("J34").Value = 2250
$ws.Range("L34").Value = 6750
$ws.Range("N34").Value = -6918
$ws.Range("H55").Value = 518.1818
$ws.Range("J55").Value = 600
$ws.Range("L55").Value = 1800
$ws.Range("N55").Value = -2154
$ws.Range("H107").Value = 610.4
$ws.Range("I107").Value = 279.2
$ws.Range("J107").Value = 657.7143
$ws.Range("K107").Value = 837.5999999999999
$ws.Range("L107").Value = 1973.1429
$ws.Range("M107").Value = 1082.4
$ws.Range("N107").Value = -5813.1429
$ws.Range("H118").Value = 2850
$ws.Range("I118").Value = 3966.6667
$ws.Range("J118").Value = 2592.3076
$ws.Range("K118").Value = 11900.0001
$ws.Range("L118").Value = 7776.9228
$ws.Range("M118").Value = -10657.0001
$ws.Range("N118").Value = -10262.9228
$ws.Range("H125").Value = 4296.25
$ws.Range("I125").Value = 3082.5
$ws.Range("J125").Value = 4498.5415
$ws.Range("K125").Value = 9247.5
$ws.Range("L125").Value = 13495.6245
$ws.Range("M125").Value = -4327.5
$ws.Range("N125").Value = -23335.6245
$ws.Range("H131").Value = 2690.6553
$ws.Range("I131").Value = 836.5833
$ws.Range("J131").Value = 3174.3262
$ws.Range("K131").Value = 2509.7499
$ws.Range("L131").Value = 9522.9786
$ws.Range("M131").Value = 2530.2501
$ws.Range("N131").Value = -19602.9786

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2364.2856
$ws.Range("I80").Value = 2306.4375
$ws.Range("J80").Value = 2441.4167
$ws.Range("K80").Value = 2306.4375
$ws.Range("L80").Value = 2441.4167
$ws.Range("M80").Value = -1308.4375
$ws.Range("N80").Value = -4437.4167
$ws.Range("H83").Value = 2364.2856
$ws.Range("I83").Value = 2306.4375
$ws.Range("J83").Value = 2441.4167
$ws.Range("K83").Value = 11532.1875
$ws.Range("L83").Value = 12207.0835
$ws.Range("M83").Value = -6540.1875
$ws.Range("N83").Value = -22191.0835
$ws.Range("H102").Value = 2420.5
$ws.Range("I102").Value = 2595.276
$ws.Range("K102").Value = 2595.276
$ws.Range("M102").Value = -973.2759999999998
$ws.Range("H126").Value = 2916.6155
$ws.Range("I126").Value = 2149.4443
$ws.Range("J126").Value = 3322.7646
$ws.Range("K126").Value = 6448.3329
$ws.Range("L126").Value = 9968.293799999999
$ws.Range("M126").Value = -3978.3329
$ws.Range("N126").Value = -14908.2938
$ws.Range("H132").Value = 3022.25
$ws.Range("I132").Value = 3255.5088
$ws.Range("J132").Value = 2135.8667
$ws.Range("K132").Value = 9766.526400000001
$ws.Range("L132").Value = 6407.6001
$ws.Range("M132").Value = -7236.526400000001
$ws.Range("N132").Value = -11467.6001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3289.0454
$ws.Range("I40").Value = 3350.4736
$ws.Range("J40").Value = 2900
$ws.Range("K40").Value = 3350.4736
$ws.Range("L40").Value = 2900
$ws.Range("M40").Value = -3214.4736
$ws.Range("N40").Value = -3172
$ws.Range("H68").Value = 2398.9285
$ws.Range("I68").Value = 1985.5
$ws.Range("J68").Value = 2950.1667
$ws.Range("K68").Value = 1985.5
$ws.Range("L68").Value = 2950.1667
$ws.Range("M68").Value = -1236.5
$ws.Range("N68").Value = -4448.1667
$ws.Range("H71").Value = 2398.9285
$ws.Range("I71").Value = 1985.5
$ws.Range("J71").Value = 2950.1667
$ws.Range("K71").Value = 9927.5
$ws.Range("L71").Value = 14750.8335
$ws.Range("M71").Value = -6183.5
$ws.Range("N71").Value = -22238.8335
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H136").Value = 5294.8237
$ws.Range("I136").Value = 2226.258
$ws.Range("J136").Value = 37003.332
$ws.Range("K136").Value = 6678.773999999999
$ws.Range("L136").Value = 111009.996
$ws.Range("M136").Value = -4128.773999999999
$ws.Range("N136").Value = -116109.996

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3614.7144
$ws.Range("I122").Value = 3524.353
$ws.Range("J122").Value = 3998.75
$ws.Range("K122").Value = 10573.059
$ws.Range("L122").Value = 11996.25
$ws.Range("M122").Value = -8123.059000000001
$ws.Range("N122").Value = -16896.25
$ws.Range("H132").Value = 1770.2941
$ws.Range("I132").Value = 1338.8292
$ws.Range("K132").Value = 4016.487599999999
$ws.Range("M132").Value = -1486.487599999999
$ws.Range("H136").Value = 1224.5857
$ws.Range("I136").Value = 680.12244
$ws.Range("J136").Value = 2495
$ws.Range("K136").Value = 2040.36732
$ws.Range("L136").Value = 7485
$ws.Range("M136").Value = 509.6326800000002
$ws.Range("N136").Value = -12585
